$d = $word.ActiveDocument

$newXml = '<w:p><w:proofErr w:type="gramStart"/><w:r><w:t>Tổng  quan</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> dự án:</w:t></w:r></w:p><w:p><w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="23"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">Destop Watch là ứng dụng đồng hồ cá nhân bao gồm các chức năng: Digital Clock, Countdown, Stop Watch, Alarm. Đặc biệt với chức năng hẹn giờ có 2 chế </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="23"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>độ :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="23"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve"> Normal and Special .Với chế độ Special , đồng hộ sẽ đưa ra các câu hỏi và đồng hồ báo thức sẽ liên tục kêu và chỉ tắt khi người dùng nhập câu trả lời đúng điều đó đem lại sự tỉnh táo cho người dùng.</w:t></w:r></w:p><w:p><w:r><w:t>Ghi chú:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Sử dụng Java Swing để phát triển ứng </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>dụng ,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> sử dụng Github/Git để quản lí dự án và quản  lí mã nguồn.</w:t></w:r></w:p>'

$r = $d.Content
$r.InsertXML($newXml)
